$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'64.343.45"
$ws.Cells.Item(2, 5).Value = "  -2.35%  "
$ws.Cells.Item(3, 4).Value = "'3.174.67"
$ws.Cells.Item(3, 5).Value = "  -7.77%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).Value = "'563.03"
$ws.Cells.Item(5, 5).Value = "  -3.70%  "
$ws.Cells.Item(6, 4).Value = "'171.60"
$ws.Cells.Item(6, 5).Value = "  -1.23%  "
$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(7, 4).Value = "'0.611"
$ws.Cells.Item(7, 5).Value = "  +2.01%  "
$ws.Cells.Item(8, 2).Value = "USDC"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$ws.Cells.Item(9, 4).Value = "'3.173.07"
$ws.Cells.Item(9, 5).Value = "  -7.73%  "
$ws.Cells.Item(11, 4).Value = "'6.61"
$ws.Cells.Item(11, 5).Value = "  -4.51%  "
$ws.Cells.Item(12, 5).Value = "  -3.52%  "
$ws.Cells.Item(13, 4).Value = "'3.729.24"
$ws.Cells.Item(13, 5).Value = "  -7.65%  "
$ws.Cells.Item(14, 5).Value = "  +1.37%  "
$ws.Cells.Item(15, 4).Value = "'27.54"
$ws.Cells.Item(15, 5).Value = "  -4.61%  "
$ws.Cells.Item(16, 4).Value = "'64.318.78"
$ws.Cells.Item(16, 5).Value = "  -2.37%  "
$ws.Cells.Item(17, 5).Value = "  -4.66%  "
$ws.Cells.Item(18, 4).Value = "'3.178.15"
$ws.Cells.Item(18, 5).Value = "  -7.87%  "
$ws.Cells.Item(19, 5).Value = "  -3.93%  "
$ws.Cells.Item(20, 4).Value = "'13.10"
$ws.Cells.Item(20, 5).Value = "  -4.77%  "
$ws.Cells.Item(21, 4).Value = "'354.12"
$ws.Cells.Item(21, 5).Value = "  -4.38%  "
$ws.Cells.Item(22, 4).Value = "'7.20"
$ws.Cells.Item(22, 5).Value = "  -5.08%  "
$ws.Cells.Item(23, 5).Value = "  +0.22%  "
$ws.Cells.Item(24, 4).Value = "'69.18"
$ws.Cells.Item(24, 5).Value = "  -4.15%  "
$ws.Cells.Item(25, 5).Value = "  -4.33%  "
$ws.Cells.Item(26, 5).Value = "  -2.65%  "
$ws.Cells.Item(27, 4).Value = "'9.74"
$ws.Cells.Item(27, 5).Value = "  +0.37%  "
$ws.Cells.Item(28, 5).Value = "  -1.39%  "
$ws.Cells.Item(29, 4).Value = "'0.997"
$ws.Cells.Item(29, 5).Value = "  -0.28%  "
$ws.Cells.Item(30, 4).Value = "'5.68"
$ws.Cells.Item(30, 5).Value = "  -1.47%  "
$ws.Cells.Item(31, 4).Value = "'0.998"
$ws.Cells.Item(31, 5).Value = "  -0.10%  "
$ws.Cells.Item(32, 5).Value = "  -3.82%  "
$ws.Cells.Item(33, 4).Value = "'22.12"
$ws.Cells.Item(33, 5).Value = "  -6.20%  "
$ws.Cells.Item(34, 5).Value = "  -4.58%  "
$ws.Cells.Item(35, 5).Value = "  -5.79%  "
$ws.Cells.Item(36, 5).Value = "  -5.24%  "
$ws.Cells.Item(37, 4).Value = "'154.77"
$ws.Cells.Item(37, 5).Value = "  -3.83%  "
$ws.Cells.Item(38, 4).Value = "'0.811"
$ws.Cells.Item(38, 5).Value = "  -7.69%  "
$ws.Cells.Item(39, 4).Value = "'25.97"
$ws.Cells.Item(39, 5).Value = "  -8.52%  "
$ws.Cells.Item(40, 4).Value = "'2.52"
$ws.Cells.Item(40, 5).Value = "  -3.93%  "
$ws.Cells.Item(41, 5).Value = "  -4.79%  "
$ws.Cells.Item(42, 4).Value = "'2.648.12"
$ws.Cells.Item(42, 5).Value = "  -4.43%  "
$ws.Cells.Item(43, 5).Value = "  -6.09%  "
$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(44, 4).Value = "'6.06"
$ws.Cells.Item(44, 5).Value = "  -5.81%  "
$ws.Cells.Item(45, 2).Value = "Bittensor"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(45, 4).Value = "'334.46"
$ws.Cells.Item(45, 5).Value = "  +3.21%  "
$ws.Cells.Item(46, 5).Value = "  -3.93%  "
$ws.Cells.Item(47, 4).Value = "'38.69"
$ws.Cells.Item(47, 5).Value = "  -3.11%  "
$ws.Cells.Item(48, 4).Value = "'23.85"
$ws.Cells.Item(48, 5).Value = "  -3.41%  "
$ws.Cells.Item(49, 5).Value = "  -6.01%  "
$ws.Cells.Item(50, 4).Value = "'0.102"
$ws.Cells.Item(50, 5).Value = "  -0.82%  "
$ws.Cells.Item(51, 4).Value = "'1.00"
$ws.Cells.Item(51, 5).Value = "  +0.04%  "
